# Update worksheet with newly recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 6.101885666666667
$ws.Range("N2").Value = 18.305657
$ws.Range("O2").Value = 0.1093737608697887
$ws.Range("P2").Value = 0.1093737608697887
$ws.Range("Q2").Value = 1.610824593372
$ws.Range("R2").Value = 14.497421340348
$ws.Range("S2").Value = 0.1093737608697887
$ws.Range("T2").Value = 0.1093737608697887

# Row 3
$ws.Range("N3").Value = 87.53628900000001
$ws.Range("O3").Value = 0.5230171820937495
$ws.Range("P3").Value = 0.5230171820937495
$ws.Range("Q3").Value = 7.702843286844001
$ws.Range("R3").Value = 69.32558958159601
$ws.Range("S3").Value = 0.5230171820937495
$ws.Range("T3").Value = 0.5230171820937495

# Row 4
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.146644
$ws.Range("N4").Value = 0.439932
$ws.Range("O4").Value = 0.002628532664354407
$ws.Range("P4").Value = 0.002628532664354407
$ws.Range("Q4").Value = 0.038712256272
$ws.Range("R4").Value = 0.348410306448
$ws.Range("S4").Value = 0.002628532664354407
$ws.Range("T4").Value = 0.002628532664354407

# Row 5
$ws.Range("M5").Value = 15.02284966666667
$ws.Range("N5").Value = 45.068549
$ws.Range("O5").Value = 0.2692783275177917
$ws.Range("P5").Value = 0.2692783275177917
$ws.Range("Q5").Value = 3.965852037804
$ws.Range("R5").Value = 35.692668340236
$ws.Range("S5").Value = 0.2692783275177917
$ws.Range("T5").Value = 0.2692783275177917

# Row 6
$ws.Range("M6").Value = 5.288900666666667
$ws.Range("N6").Value = 15.866702
$ws.Range("O6").Value = 0.09480134312252211
$ws.Range("P6").Value = 0.09480134312252211
$ws.Range("Q6").Value = 1.396206309192
$ws.Range("R6").Value = 12.565856782728
$ws.Range("S6").Value = 0.09480134312252211
$ws.Range("T6").Value = 0.09480134312252211

# Row 7
$ws.Range("M7").Value = 0.050258
$ws.Range("N7").Value = 0.150774
$ws.Range("O7").Value = 0.0009008537317934847
$ws.Range("P7").Value = 0.0009008537317934848
$ws.Range("Q7").Value = 0.013267508904
$ws.Range("R7").Value = 0.119407580136
$ws.Range("S7").Value = 0.0009008537317934847
$ws.Range("T7").Value = 0.0009008537317934848
